# Auto-generated edit script applying the diff's numeric updates
# to the Pandaemonium_Profits workbook (sheets ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H20").Value = 10560
$ws.Range("J20").Value = 12900
$ws.Range("L20").Value = 12900
$ws.Range("N20").Value = -13360

$ws.Range("H35").Value = 10560
$ws.Range("J35").Value = 12900
$ws.Range("L35").Value = 12900
$ws.Range("N35").Value = -13658

$ws.Range("H100").Value = 1387.6428
$ws.Range("I100").Value = 1325.2222
$ws.Range("J100").Value = 1500
$ws.Range("K100").Value = 1325.2222
$ws.Range("L100").Value = 1500
$ws.Range("M100").Value = -784.2221999999999
$ws.Range("N100").Value = -2582

$ws.Range("H138").Value = 2927.1
$ws.Range("I138").Value = 1256.5264
$ws.Range("J138").Value = 5812.636
$ws.Range("K138").Value = 3769.5792
$ws.Range("L138").Value = 17437.908
$ws.Range("M138").Value = 1370.4208
$ws.Range("N138").Value = -27717.908

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 20181.18
$ws.Range("I32").Value = 22700.404
$ws.Range("J32").Value = 11447.866
$ws.Range("K32").Value = 22700.404
$ws.Range("L32").Value = 11447.866
$ws.Range("M32").Value = -22413.404
$ws.Range("N32").Value = -12021.866

$ws.Range("H41").Value = 7599.3335
$ws.Range("I41").Value = 1499
$ws.Range("K41").Value = 1499
$ws.Range("M41").Value = -1085

$ws.Range("H42").Value = 40000
$ws.Range("J42").Value = 40000
$ws.Range("L42").Value = 40000
$ws.Range("N42").Value = -40972

$ws.Range("H122").Value = 2607.2
$ws.Range("I122").Value = 3320
$ws.Range("J122").Value = 2250.8
$ws.Range("K122").Value = 9960
$ws.Range("L122").Value = 6752.400000000001
$ws.Range("M122").Value = -7510
$ws.Range("N122").Value = -11652.4

$ws.Range("H132").Value = 4115.357
$ws.Range("I132").Value = 1872.5264
$ws.Range("J132").Value = 8850.223
$ws.Range("K132").Value = 5617.5792
$ws.Range("L132").Value = 26550.669
$ws.Range("M132").Value = -3087.5792
$ws.Range("N132").Value = -31610.669

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 28078.205
$ws.Range("I134").Value = 2321.4482
$ws.Range("J134").Value = 102772.8
$ws.Range("K134").Value = 6964.344599999999
$ws.Range("L134").Value = 308318.4
$ws.Range("M134").Value = -4429.344599999999
$ws.Range("N134").Value = -313388.4

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H41").Value = 29000
$ws.Range("I41").Value = 0
$ws.Range("J41").Value = 29000
$ws.Range("K41").Value = 0
$ws.Range("L41").Value = 29000
$ws.Range("M41").ClearContents()
$ws.Range("N41").Value = -29856

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H41").Value = 477
$ws.Range("J41").Value = 763.5
$ws.Range("L41").Value = 2290.5
$ws.Range("N41").Value = -2966.5

$ws.Range("H121").Value = 597.4545000000001
$ws.Range("J121").Value = 711.625
$ws.Range("L121").Value = 2134.875
$ws.Range("N121").Value = -4754.875

$ws.Range("H131").Value = 1021.0571
$ws.Range("I131").Value = 861.1111
$ws.Range("J131").Value = 1044.6558
$ws.Range("K131").Value = 2583.3333
$ws.Range("L131").Value = 3133.9674
$ws.Range("M131").Value = 2456.6667
$ws.Range("N131").Value = -13213.9674

$ws.Range("H137").Value = 19658.594
$ws.Range("I137").Value = 1649.2667
$ws.Range("J137").Value = 35549.176
$ws.Range("K137").Value = 4947.800099999999
$ws.Range("L137").Value = 106647.528
$ws.Range("M137").Value = 152.1999000000005
$ws.Range("N137").Value = -116847.528

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H41").Value = 14525
$ws.Range("I41").Value = 3500
$ws.Range("J41").Value = 18200
$ws.Range("K41").Value = 3500
$ws.Range("L41").Value = 18200
$ws.Range("M41").Value = -3145
$ws.Range("N41").Value = -18910

$ws.Range("H43").Value = 8408.5
$ws.Range("I43").Value = 2017
$ws.Range("J43").Value = 14800
$ws.Range("K43").Value = 2017
$ws.Range("L43").Value = 14800
$ws.Range("M43").Value = -1866
$ws.Range("N43").Value = -15102

$ws.Range("H44").Value = 0
$ws.Range("J44").Value = 0
$ws.Range("L44").Value = 0
$ws.Range("N44").ClearContents()

$ws.Range("H105").Value = 60671
$ws.Range("J105").Value = 60671
$ws.Range("L105").Value = 60671
$ws.Range("N105").Value = -67659

$ws.Range("H132").Value = 8093.9116
$ws.Range("I132").Value = 5472.1035
$ws.Range("J132").Value = 23300.4
$ws.Range("K132").Value = 16416.3105
$ws.Range("L132").Value = 69901.20000000001
$ws.Range("M132").Value = -13886.3105
$ws.Range("N132").Value = -74961.20000000001

$ws.Range("H139").Value = 39663
$ws.Range("J139").Value = 39663
$ws.Range("L139").Value = 39663
$ws.Range("N139").Value = -49943

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 948.05554
$ws.Range("I22").Value = 900.5
$ws.Range("J22").Value = 971.8333
$ws.Range("K22").Value = 900.5
$ws.Range("L22").Value = 971.8333
$ws.Range("M22").Value = -605.5
$ws.Range("N22").Value = -1561.8333

$ws.Range("H27").Value = 948.05554
$ws.Range("I27").Value = 900.5
$ws.Range("J27").Value = 971.8333
$ws.Range("K27").Value = 900.5
$ws.Range("L27").Value = 971.8333
$ws.Range("M27").Value = -793.5
$ws.Range("N27").Value = -1185.8333

$ws.Range("H42").Value = 22400
$ws.Range("J42").Value = 22400
$ws.Range("L42").Value = 22400
$ws.Range("N42").Value = -23526

$ws.Range("H49").Value = 22400
$ws.Range("J49").Value = 22400
$ws.Range("L49").Value = 22400
$ws.Range("N49").Value = -22694

$ws.Range("H141").Value = 55150
$ws.Range("J141").Value = 55150
$ws.Range("L141").Value = 55150
$ws.Range("N141").Value = -65510

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H15").Value = 8800
$ws.Range("J15").Value = 10000
$ws.Range("L15").Value = 10000
$ws.Range("N15").Value = -10576

$ws.Range("H99").Value = 62738
$ws.Range("I99").Value = 0
$ws.Range("J99").Value = 62738
$ws.Range("K99").Value = 0
$ws.Range("L99").Value = 62738
$ws.Range("M99").ClearContents()
$ws.Range("N99").Value = -68728

$ws.Range("H140").Value = 54597
$ws.Range("J140").Value = 54597
$ws.Range("L140").Value = 54597
$ws.Range("N140").Value = -64957

$ws.Range("H141").Value = 59943
$ws.Range("J141").Value = 59943
$ws.Range("L141").Value = 59943
$ws.Range("N141").Value = -70303
